$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D: "success" header + per-row success flag ("0"/"1" as text)

# D1 header — same text + bold/border/centered style as the other headers
$ws.Range("D1").Value = "success"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# D2:D4 — values must be stored as text ("0"/"1"), not numbers, and keep
# the plain (unstyled) look of the other body cells. Formatting the cell
# as Text before assigning the value keeps Excel from coercing the
# numeric-looking string into a real number; switching back to the
# "Normal" style afterward drops the temporary number format again.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1"
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0"
$ws.Range("D4").Style = "Normal"
